$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Update the email value in B2 and B5 from ptvanh9@mailinator.com to ptvanh@mailinator.com
$ws.Range("B2").Value = "ptvanh@mailinator.com"
$ws.Range("B5").Value = "ptvanh@mailinator.com"

# Update the selected cell on the sheet to F17
$ws.Activate()
$ws.Range("F17").Select()
